# Applies the commit: inserts 3 new weekly "Espárragos" price rows for the
# Mercado Mayorista Lo Valledor de Santiago sheet, just before the old row 90,
# pushing the remaining data rows down by 3 (old row 90 -> new row 93, etc.),
# and fills the 3 freshly inserted rows with the new reporting data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 90:92 (existing rows 90.. shift down to 93..)
$ws.Range("90:92").Insert()

# ---- New row 90 ----
$ws.Cells.Item(90, 1).Value2  = 6
$ws.Cells.Item(90, 2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(90, 3).Value2  = "Metropolitana"
$ws.Cells.Item(90, 4).Value2  = 44907
$ws.Cells.Item(90, 5).Value2  = 13
$ws.Cells.Item(90, 6).Value2  = 300000000
$ws.Cells.Item(90, 7).Value2  = "Espárragos"
$ws.Cells.Item(90, 8).Value2  = "Sin especificar"
$ws.Cells.Item(90, 9).Value2  = "Primera"
$ws.Cells.Item(90, 10).Value2 = 1600
$ws.Cells.Item(90, 11).Value2 = 1800
$ws.Cells.Item(90, 12).Value2 = 1800
$ws.Cells.Item(90, 13).Value2 = 1800
$ws.Cells.Item(90, 14).Value2 = "$/kilo"
$ws.Cells.Item(90, 15).Value2 = "Provincia de Linares"
$ws.Cells.Item(90, 16).Value2 = 1800
$ws.Cells.Item(90, 17).Value2 = 1
$ws.Cells.Item(90, 18).Value2 = "Hortaliza"

# ---- New row 91 ----
$ws.Cells.Item(91, 1).Value2  = 6
$ws.Cells.Item(91, 2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(91, 3).Value2  = "Metropolitana"
$ws.Cells.Item(91, 4).Value2  = 44907
$ws.Cells.Item(91, 5).Value2  = 13
$ws.Cells.Item(91, 6).Value2  = 300000000
$ws.Cells.Item(91, 7).Value2  = "Espárragos"
$ws.Cells.Item(91, 8).Value2  = "Sin especificar"
$ws.Cells.Item(91, 9).Value2  = "Segunda"
$ws.Cells.Item(91, 10).Value2 = 1200
$ws.Cells.Item(91, 11).Value2 = 1500
$ws.Cells.Item(91, 12).Value2 = 1500
$ws.Cells.Item(91, 13).Value2 = 1500
$ws.Cells.Item(91, 14).Value2 = "$/kilo"
$ws.Cells.Item(91, 15).Value2 = "Provincia de Linares"
$ws.Cells.Item(91, 16).Value2 = 1500
$ws.Cells.Item(91, 17).Value2 = 1
$ws.Cells.Item(91, 18).Value2 = "Hortaliza"

# ---- New row 92 ----
$ws.Cells.Item(92, 1).Value2  = 6
$ws.Cells.Item(92, 2).Value2  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(92, 3).Value2  = "Metropolitana"
$ws.Cells.Item(92, 4).Value2  = 44509
$ws.Cells.Item(92, 5).Value2  = 13
$ws.Cells.Item(92, 6).Value2  = 300000000
$ws.Cells.Item(92, 7).Value2  = "Espárragos"
$ws.Cells.Item(92, 8).Value2  = "Sin especificar"
$ws.Cells.Item(92, 9).Value2  = "Tercera"
$ws.Cells.Item(92, 10).Value2 = 740
$ws.Cells.Item(92, 11).Value2 = 1200
$ws.Cells.Item(92, 12).Value2 = 1200
$ws.Cells.Item(92, 13).Value2 = 1200
$ws.Cells.Item(92, 14).Value2 = "$/kilo"
$ws.Cells.Item(92, 15).Value2 = "Provincia de Linares"
$ws.Cells.Item(92, 16).Value2 = 1200
$ws.Cells.Item(92, 17).Value2 = 1
$ws.Cells.Item(92, 18).Value2 = "Hortaliza"
